$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking price strings
# (e.g. "51.043.50") are not auto-converted to numbers by Excels type inference.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('B2').Value = 'Bitcoin'
$ws.Range('C2').Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range('D2').Value = '51.043.50'
$ws.Range('E2').Value = '  -1.56%  '

$ws.Range('B3').Value = 'Ethereum'
$ws.Range('C3').Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range('D3').Value = '2.942.32'
$ws.Range('E3').Value = '  -2.21%  '

$ws.Range('B4').Value = 'TetherUSD'
$ws.Range('C4').Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = '376.85'
$ws.Range('E5').Value = '  -0.96%  '

$ws.Range('B6').Value = 'Solana'
$ws.Range('C6').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D6').Value = '102.48'
$ws.Range('E6').Value = '  -3.44%  '

$ws.Range('B7').Value = 'XRP'
$ws.Range('C7').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D7').Value = '0.538'
$ws.Range('E7').Value = '  -1.68%  '

$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.03%  '

$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').Value = '0.586'
$ws.Range('E9').Value = '  -3.00%  '

$ws.Range('B10').Value = 'Avalanche'
$ws.Range('C10').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D10').Value = '36.55'
$ws.Range('E10').Value = '  -3.92%  '

$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').Value = '0.138'
$ws.Range('E11').Value = '  -0.84%  '

$ws.Range('B12').Value = 'Dogecoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D12').Value = '0.0839'
$ws.Range('E12').Value = '  -1.32%  '

$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '3.399.98'
$ws.Range('E13').Value = '  -1.90%  '

$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = '18.01'
$ws.Range('E14').Value = '  -4.62%  '

$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '7.38'
$ws.Range('E15').Value = '  -2.49%  '

$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.929.87'
$ws.Range('E16').Value = '  -2.06%  '

$ws.Range('B17').Value = 'Polygon'
$ws.Range('C17').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D17').Value = '0.981'
$ws.Range('E17').Value = '  +1.19%  '

$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '50.941.62'
$ws.Range('E18').Value = '  -1.87%  '

$ws.Range('B19').Value = 'ImmutableX'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D19').Value = '3.18'
$ws.Range('E19').Value = '  -9.78%  '

$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').Value = '7.13'
$ws.Range('E20').Value = '  -4.77%  '

$ws.Range('B21').Value = 'InternetComputer(DFINITY)'
$ws.Range('C21').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D21').Value = '12.53'
$ws.Range('E21').Value = '  -5.49%  '

$ws.Range('B22').Value = 'ShibaInu'
$ws.Range('C22').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D22').Value = '0.0₃0952'
$ws.Range('E22').Value = '  -1.28%  '

$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').Value = '68.38'
$ws.Range('E23').Value = '  -0.93%  '

$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').Value = '262.21'
$ws.Range('E24').Value = '  -1.09%  '

$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').Value = '2.89'
$ws.Range('E25').Value = '  +2.85%  '

$ws.Range('B26').Value = 'Filecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D26').Value = '8.24'
$ws.Range('E26').Value = '  +8.89%  '

$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').Value = '7.64'
$ws.Range('E27').Value = '  +2.29%  '

$ws.Range('B28').Value = 'Hedera'
$ws.Range('C28').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D28').Value = '0.116'
$ws.Range('E28').Value = '  +9.41%  '

$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').Value = '0.168'
$ws.Range('E29').Value = '  -2.87%  '

$ws.Range('B30').Value = 'Dai'
$ws.Range('C30').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.03%  '

$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '25.62'
$ws.Range('E31').Value = '  -2.55%  '

$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D32').Value = '9.79'
$ws.Range('E32').Value = '  -2.21%  '

$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D33').Value = '34.20'
$ws.Range('E33').Value = '  -2.33%  '

$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').Value = '50.58'
$ws.Range('E34').Value = '  -2.11%  '

$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D35').Value = '0.0454'
$ws.Range('E35').Value = '  +2.81%  '

$ws.Range('B36').Value = 'Toncoin'
$ws.Range('C36').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D36').Value = '2.05'
$ws.Range('E36').Value = '  -1.66%  '

$ws.Range('B37').Value = 'FirstDigitalUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  -0.09%  '

$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').Value = '2.98'
$ws.Range('E38').Value = '  -4.76%  '

$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = '2.59'
$ws.Range('E39').Value = '  -3.95%  '

$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D40').Value = '16.61'
$ws.Range('E40').Value = '  -5.48%  '

$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').Value = '0.115'
$ws.Range('E41').Value = '  -1.42%  '

$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').Value = '1.78'
$ws.Range('E42').Value = '  -5.50%  '

$ws.Range('B43').Value = 'Monero'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D43').Value = '120.94'
$ws.Range('E43').Value = '  -2.94%  '

$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '21.36'
$ws.Range('E44').Value = '  -4.49%  '

$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').Value = '2.05'
$ws.Range('E45').Value = '  -1.90%  '

$ws.Range('B46').Value = 'TheGraph'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D46').Value = '0.274'
$ws.Range('E46').Value = '  -5.18%  '

$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value = '2.38'
$ws.Range('E47').Value = '  -1.52%  '

$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '3.23'
$ws.Range('E48').Value = '  -1.99%  '

$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '2.006.14'
$ws.Range('E49').Value = '  -2.29%  '

$ws.Range('B50').Value = 'BEAM'
$ws.Range('C50').Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range('D50').Value = '0.0345'
$ws.Range('E50').Value = '  +1.09%  '

$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D51').Value = '0.482'
$ws.Range('E51').Value = '  +12.28%  '

# Restore the default (Normal) cell style on column D so no stray
# number-format style is left referenced on these cells.
$ws.Range('D2:D51').Style = 'Normal'
